$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows above row 78, pushing the existing rows 78-89
# down to 82-93 (and carrying the date-format style from row 78's "D"
# column along with them).
$ws.Rows("78:81").Insert()

# New weekly entries (Provincia del Elquí, same Mercado/Producto block as
# the surrounding rows) for the 2021-11-04 reporting date.
$newRows = @(
    @{ Row=78; D=44504; L='Especial'; M=50;  N=3500;  O=3500;  P=3500;  Q='$/kilo (en caja de 15 kilos)'; R='Provincia del Elquí'; S=3500; T=1 },
    @{ Row=79; D=44504; L='Primera';  M=160; N=3000;  O=3000;  P=3000;  Q='$/kilo (en caja de 15 kilos)'; R='Provincia del Elquí'; S=3000; T=1 },
    @{ Row=80; D=44504; L='Segunda';  M=80;  N=20000; O=20000; P=20000; Q='$/bandeja 8 kilos';            R='Provincia del Elquí'; S=2500; T=8 },
    @{ Row=81; D=44504; L='Segunda';  M=80;  N=2500;  O=2500;  P=2500;  Q='$/kilo (en caja de 15 kilos)'; R='Provincia del Elquí'; S=2500; T=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 10
    $ws.Cells.Item($row, 2).Value = 'Vega Modelo de Temuco'
    $ws.Cells.Item($row, 3).Value = 'La Araucanía'
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 9
    $ws.Cells.Item($row, 6).Value = 'Fruta'
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = 'Otros'
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = 'Chirimoya'
    $ws.Cells.Item($row, 11).Value = 'Cultivar IV Región'
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
